$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.62
$ws.Range("S2").Value = 3.1
$ws.Range("U2").Value = 1.73
$ws.Range("AD2").Value = 34

$ws.Range("J3").Value = 7.8
$ws.Range("K3").Value = 9
$ws.Range("R3").Value = 2.02
$ws.Range("T3").Value = 1.63
$ws.Range("AN3").Value = 3.2

$ws.Range("AA4").Value = 30
$ws.Range("AO4").Value = 25

$ws.Range("R5").Value = 1.26
$ws.Range("X5").Value = 11

$ws.Range("Z6").Value = 16
$ws.Range("AM6").Value = 120
